$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: new bug-tracker entry -------------------------------------------------
# Bug ID
$ws.Range("A6").Value = 5

# Issue
$ws.Range("B6").Value = "First question is not formatted correctly."

# Date - copy the date formatting from the row above first so we reuse the
# existing date-formatted cell style instead of creating a brand new one,
# then assign the date itself.
$ws.Range("C5").Copy() | Out-Null
$ws.Range("C6").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("C6").Value = (Get-Date -Year 2024 -Month 7 -Day 8).Date

# Description
$ws.Range("D6").Value = "First question is not formatted correctly."

# Status
$ws.Range("E6").Value = "Ongoing"

# The new row's content needs two lines to display, so the row grows taller.
$ws.Rows("6").RowHeight = 29

# --- Selection / active cell -------------------------------------------------------
$ws.Range("G6").Select() | Out-Null
